# Add a new worksheet "t9_ML_2012" at the end of the workbook (after "t8_ML_1612")
# with data for 3 students, and update the selection/active-tab state so the new
# sheet becomes the tab shown when the workbook is opened (matching the
# "tabSelected" / "activeTab" bookkeeping in the target workbook).

$wb = $excel.ActiveWorkbook

# --- Update selection on the previously-last sheet (t8_ML_1612) first, before
# --- activating/adding the new sheet, so it ends up no longer the selected tab.
$ws8 = $wb.Worksheets.Item("t8_ML_1612")
$ws8.Range("F2:F7").Style = "Normal"
$ws8.Range("H22").Select() | Out-Null

# --- Insert the new sheet right after the current last sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "t9_ML_2012"

# Header row
$ws.Range("A1").Value = "Nom de famille"
$ws.Range("B1").Value = "Prénom"
$ws.Range("C1").Value = "Clé"
$ws.Range("D1").Value = "Adresse de courriel"
$ws.Range("E1").Value = "Durée"
$ws.Range("F1").Value = "Note/20,00"
$ws.Range("G1").Value = "Q. 1 /2,00"
$ws.Range("H1").Value = "Q. 2 /2,00"
$ws.Range("I1").Value = "Q. 3 /2,50"
$ws.Range("J1").Value = "Q. 4 /2,50"
$ws.Range("K1").Value = "Q. 5 /3,00"
$ws.Range("L1").Value = "Q. 6 /3,00"
$ws.Range("M1").Value = "Q. 7 /2,00"
$ws.Range("N1").Value = "Q. 8 /3,00"

# Row 2 - Bertille PECOURT
$ws.Range("A2").Value = "PECOURT"
$ws.Range("B2").Value = "Bertille"
$ws.Range("C2").Formula = "=A2&B2"
$ws.Range("D2").Value = "bertille.pecourt@etu.unilasalle.fr"
$ws.Range("E2").Value = "7 min 28 s"
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0

# Row 3 - Jules CHIEN-CHOW-CHINE
$ws.Range("A3").Value = "CHIEN-CHOW-CHINE"
$ws.Range("B3").Value = "Jules"
$ws.Range("D3").Value = "jules.chien-chow-chine@etu.unilasalle.fr"
$ws.Range("E3").Value = "10 min 8 s"
$ws.Range("F3").Value = 6
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 3

# Row 4 - Hugo BOUTILLIER
$ws.Range("A4").Value = "BOUTILLIER"
$ws.Range("B4").Value = "Hugo"
$ws.Range("D4").Value = "hugo.boutillier@etu.unilasalle.fr"
$ws.Range("E4").Value = "7 min 33 s"
$ws.Range("F4").Value = 8.5
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 2.5
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 3

# C3:C4 share a formula (A3&B3) - set as one range write to produce a shared formula
$ws.Range("C3:C4").Formula = "=A3&B3"

# Selection on the new sheet (becomes the active / last-selected tab)
$ws.Range("L16").Select() | Out-Null
